$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF headers using same style as existing headers (style index 1)
# Copy the H1 header cell's formatting onto I1/J1, then overwrite the value/text.
$ws.Cells.Item(1, 8).Copy($ws.Cells.Item(1, 9))
$ws.Cells.Item(1, 8).Copy($ws.Cells.Item(1, 10))

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-48
$data = @(
    @(7, 7),
    @(5, 7),
    @(8, 8),
    @(8, 8),
    @(4, 5),
    @(7, 8),
    @(9, 9),
    @(6, 6),
    @(7, 8),
    @(8, 9),
    @(7, 7),
    @(8, 9),
    @(8, 9),
    @(5, 6),
    @(4, 6),
    @(7, 8),
    @(5, 6),
    @(6, 7),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(8, 10),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(1, 3),
    @(5, 6),
    @(6, 6),
    @(5, 6),
    @(6, 6),
    @(3, 3),
    @(3, 4),
    @(7, 7),
    @(8, 8),
    @(4, 4),
    @(3, 3)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
